$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scenario run modeling of renovation changes: update calibrated dwelling
# count projections for column B (rows 303-452, years 1901-2050).

$ws.Range("B303").Value = 1562.576980323525
$ws.Range("B304").Value = 1587.38525962832
$ws.Range("B305").Value = 1612.623564947919
$ws.Range("B306").Value = 1638.297748936454
$ws.Range("B307").Value = 1664.413561350409
$ws.Range("B308").Value = 1690.976633843063
$ws.Range("B309").Value = 1717.99246416102
$ws.Range("B310").Value = 1745.466399747083
$ws.Range("B311").Value = 1773.403620756999
$ws.Range("B312").Value = 1801.809122503792
$ws.Range("B313").Value = 1830.68769734128
$ws.Range("B314").Value = 1860.043916002599
$ws.Range("B315").Value = 1889.882108416085
$ws.Range("B316").Value = 1920.20634402
$ws.Range("B317").Value = 1951.020411602156
$ws.Range("B318").Value = 1982.327798697268
$ws.Range("B319").Value = 2014.13167057277
$ws.Range("B320").Value = 2046.434848843066
$ws.Range("B321").Value = 2079.239789755764
$ws.Range("B322").Value = 2112.548562192103
$ws.Range("B323").Value = 2146.362825437073
$ws.Range("B324").Value = 2180.683806769377
$ws.Range("B325").Value = 2215.51227893648
$ws.Range("B326").Value = 2250.848537573061
$ws.Range("B327").Value = 2286.692378635208
$ws.Range("B328").Value = 2323.043075922918
$ws.Range("B329").Value = 2359.899358767143
$ws.Range("B330").Value = 2397.259389964813
$ws.Range("B331").Value = 2435.120744048032
$ws.Range("B332").Value = 2473.480385978243
$ws.Range("B333").Value = 2512.334650358807
$ws.Range("B334").Value = 2551.679221266879
$ws.Range("B335").Value = 2591.50911280471
$ws.Range("B336").Value = 2631.818650476983
$ws.Range("B337").Value = 2672.601453503303
$ws.Range("B338").Value = 2713.850418175942
$ws.Range("B339").Value = 2755.557702376225
$ws.Range("B340").Value = 2797.714711368685
$ws.Range("B341").Value = 2840.312084985015
$ws.Range("B342").Value = 2883.339686318621
$ws.Range("B343").Value = 2926.786592048765
$ws.Range("B344").Value = 2970.641084510163
$ws.Range("B345").Value = 3014.890645628846
$ws.Range("B346").Value = 3059.521952835146
$ws.Range("B347").Value = 3104.520877075838
$ws.Range("B348").Value = 3149.872483027312
$ws.Range("B349").Value = 3195.561031627502
$ws.Range("B350").Value = 3241.569985020744
$ws.Range("B351").Value = 3287.882014023539
$ws.Range("B352").Value = 3334.479008197065
$ws.Range("B353").Value = 3381.342088618478
$ws.Range("B354").Value = 3428.451623426741
$ws.Range("B355").Value = 3475.787246218047
$ws.Range("B356").Value = 3523.327877351266
$ws.Range("B357").Value = 3571.051748218246
$ws.Range("B358").Value = 4109.747156813833
$ws.Range("B359").Value = 4164.282526845119
$ws.Range("B360").Value = 4218.947457957644
$ws.Range("B361").Value = 4273.7144440119
$ws.Range("B362").Value = 4328.55538605475
$ws.Range("B363").Value = 4383.441640660688
$ws.Range("B364").Value = 4438.344071709143
$ws.Range("B365").Value = 4493.233105538183
$ws.Range("B366").Value = 4548.078789413124
$ws.Range("B367").Value = 4602.850853231757
$ws.Range("B368").Value = 4657.518774357688
$ws.Range("B369").Value = 4712.051845479154
$ws.Range("B370").Value = 4766.419245350775
$ws.Range("B371").Value = 4820.590112278523
$ws.Range("B372").Value = 4874.53362018392
$ws.Range("B373").Value = 4168.589554182335
$ws.Range("B374").Value = 4213.755879528547
$ws.Range("B375").Value = 4258.652520523429
$ws.Range("B376").Value = 4303.254204546865
$ws.Range("B377").Value = 4347.536037806946
$ws.Range("B378").Value = 4391.473579182299
$ws.Range("B379").Value = 4435.042914950384
$ws.Range("B380").Value = 4478.220734178165
$ws.Range("B381").Value = 4520.984404545542
$ws.Range("B382").Value = 4563.312048363632
$ws.Range("B383").Value = 3363.071004756454
$ws.Range("B384").Value = 3393.299728800471
$ws.Range("B385").Value = 3423.165961186296
$ws.Range("B386").Value = 3452.656481248294
$ws.Range("B387").Value = 3481.758892890191
$ws.Range("B388").Value = 3510.461677190985
$ws.Range("B389").Value = 3538.754243890691
$ws.Range("B390").Value = 3566.626981577801
$ws.Range("B391").Value = 3594.07130640003
$ws.Range("B392").Value = 3621.079709127204
$ws.Range("B393").Value = 5920.129258485421
$ws.Range("B394").Value = 5962.519671805097
$ws.Range("B395").Value = 6004.177203691216
$ws.Range("B396").Value = 6045.097343877086
$ws.Range("B397").Value = 6085.277658550986
$ws.Range("B398").Value = 6124.717845129083
$ws.Range("B399").Value = 6163.419782577375
$ws.Range("B400").Value = 6201.387577048393
$ws.Range("B401").Value = 6238.627602632813
$ws.Range("B402").Value = 6275.148536999833
$ws.Range("B403").Value = 16384.06406011116
$ws.Range("B404").Value = 16475.2352635508
$ws.Range("B405").Value = 16564.64379723768
$ws.Range("B406").Value = 16652.33675348625
$ws.Range("B407").Value = 16738.36745114018
$ws.Range("B408").Value = 16822.79543806105
$ws.Range("B409").Value = 16905.68647696257
$ws.Range("B410").Value = 16987.11251405175
$ws.Range("B411").Value = 17067.15163007672
$ws.Range("B412").Value = 17145.88797322696
$ws.Range("B413").Value = 16405.64554180798
$ws.Range("B414").Value = 16478.42480507048
$ws.Range("B415").Value = 16550.23737802795
$ws.Range("B416").Value = 16621.18553030348
$ws.Range("B417").Value = 16691.37671580532
$ws.Range("B418").Value = 16760.92339505756
$ws.Range("B419").Value = 16829.94283674865
$ws.Range("B420").Value = 16898.55689793379
$ws.Range("B421").Value = 16966.89178233877
$ws.Range("B422").Value = 17035.07777617706
$ws.Range("B423").Value = 17103.24896095234
$ws.Range("B424").Value = 17171.54290264174
$ws.Range("B425").Value = 17240.10031677472
$ws.Range("B426").Value = 17309.06470883364
$ws.Range("B427").Value = 17378.58198953933
$ws.Range("B428").Value = 17448.80006457694
$ws.Range("B429").Value = 17519.86839837265
$ws.Range("B430").Value = 17591.93755165499
$ws.Range("B431").Value = 17665.1586925798
$ws.Range("B432").Value = 17739.68308132217
$ws.Range("B433").Value = 17815.66152819534
$ws.Range("B434").Value = 17893.24382544903
$ws.Range("B435").Value = 17972.57815314291
$ws.Range("B436").Value = 18053.81045959192
$ws.Range("B437").Value = 18137.08381720096
$ws.Range("B438").Value = 18222.53775462758
$ws.Range("B439").Value = 18310.30756658411
$ws.Range("B440").Value = 18400.52360275779
$ws.Range("B441").Value = 18493.31053773097
$ws.Range("B442").Value = 18588.78662403392
$ws.Range("B443").Value = 18687.06293081792
$ws.Range("B444").Value = 18788.24257102259
$ws.Range("B445").Value = 18892.41992023675
$ws.Range("B446").Value = 18999.67983085584
$ws.Range("B447").Value = 19110.0968455128
$ws.Range("B448").Value = 19223.73441415162
$ws.Range("B449").Value = 19340.6441195282
$ws.Range("B450").Value = 19460.86491620301
$ws.Range("B451").Value = 19584.42238860781
$ws.Range("B452").Value = 19711.32803397847
